# "colors added to C version"
#
# Applies the tracked changes to slide 1 ("2024 splash screen"):
#   - moves the logo picture down a bit
#   - fixes the capitalisation of "physiomon" -> "PhysioMon"
#   - repositions / resizes the author textbox and normalizes its run
#     font sizes back to the (inherited) default size
#
# NB: the revision-history stream (ppt/changesInfos/*) and the
# auto-generated "datetimeFigureOut" field cached inside the notes
# master are produced by real PowerPoint's collaboration/autosave
# machinery and are not reachable from the scripted object model, so
# they are intentionally left untouched here.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Picture 2 (id=3): nudge down -------------------------------------
$pic = $s.Shapes.Item(1)
$pic.Top = 88.5

# --- TextBox 10 (id=11): "physiomon" -> "PhysioMon" --------------------
$title = $s.Shapes.Item(2)
$title.TextFrame.TextRange.Paragraphs(1).Runs(1).Text = "PhysioMon"

# --- TextBox 11 (id=12): reposition/resize + drop explicit sz=14pt -----
$credit = $s.Shapes.Item(5)

# normalize every run's font size back to the inherited default (18pt)
$credit.TextFrame.TextRange.Font.Size = 18

# move/resize the box (literals chosen so the emitted EMU matches exactly)
$credit.Left = 446.52734375
$credit.Top = 250.5
$credit.Width = 264.0000915527344
$credit.Height = 72.70315551757812
